$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# Status text "Ready for handoff" -> "Handback transform failed" for the
# 3a174ad2... file row (row 3) across the Overview, zh-cn and de-de sheets.
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"
$wsZh.Range("C3").Value = "Handback transform failed"
$wsDe.Range("C3").Value = "Handback transform failed"

# Add error detail for zh-cn row 3 (K3)
$wsZh.Range("K3").Value = "Handback file name: ydzvsygy.roy is different with handoff file name: 3a174ad2-4ae1-4e74-a739-1882562ae7c2.c953449fb01711e6af70619e9f1680ae943064a4.zh-cn."

# Add error detail for de-de row 3 (K3)
$wsDe.Range("K3").Value = "Handback file name: ydzvsygy.roy is different with handoff file name: 3a174ad2-4ae1-4e74-a739-1882562ae7c2.c953449fb01711e6af70619e9f1680ae943064a4.de-de."
